# Gestionar Solicitudes de Amistad - corrige tildes/erratas y separa los
# runs en torno a las palabras corregidas (como deja Word al aceptar una
# sugerencia del corrector ortografico).

$d = $word.ActiveDocument

function Split-RunAt($findRange, $searchText) {
    # Busca $searchText dentro de $findRange (un Range "vivo" anclado al
    # documento) y fuerza que quede en su propio run aplicando y
    # revirtiendo una propiedad de caracter -- Word fusiona runs
    # adyacentes con el mismo formato, así que el "round trip" es lo que
    # produce el split real en tres <w:r>.
    $ok = $findRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "No se encontro '$searchText'"
    }
    $findRange.Font.Bold = 1
    $findRange.Font.Bold = 0
    return $findRange
}

function Get-CellDocRange($cell) {
    # Devuelve un Range anclado al documento (no a la celda) con los
    # mismos limites que la celda -- buscar directamente sobre
    # cell.Range dentro de una tabla no deja "prender" Find/formato.
    $s = $cell.Range.Start
    $e = $cell.Range.End
    return $d.Range($s, $e)
}

# ---------------------------------------------------------------------
# 1) "Usuario logueado selecciona..." -> separa "logueado"
# ---------------------------------------------------------------------
$rng = $d.Content
Split-RunAt $rng "logueado" | Out-Null

# ---------------------------------------------------------------------
# 2) Tabla de casos de prueba
# ---------------------------------------------------------------------
$t = $d.Tables(1)

# Fila 2, Col 1: "Click en Solicitudes" -> separa "Click"
$cell = $t.Cell(2, 1)
$rng = Get-CellDocRange $cell
Split-RunAt $rng "Click" | Out-Null

# Fila 2, Col 2: "Se mostro..." -> "Se mostró..." separa "mostro"/"mostró"
$cell = $t.Cell(2, 2)
$rng = Get-CellDocRange $cell
$rng2 = Split-RunAt $rng "mostro"
$rng2.Text = "mostró"

# Fila 3, Col 2: "Se borro la solicitud de la lista. Se agrego al usuario
# como amigo" -> corrige "borro"->"borró" y "agrego"->"agrego" con acento
$cell = $t.Cell(3, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End

$rng = $d.Range($cellStart, $cellEnd)
$rngBorro = Split-RunAt $rng "borro"
$rngBorro.Text = "borró"

$rng = $d.Range($cellStart, $cellEnd)
$rngAgrego = Split-RunAt $rng "agrego"
$rngAgrego.Text = "agregó"

# Fila 4, Col 2: "Se borro la solicitud de la lista" -> "borró" (sin bookmark)
$cell = $t.Cell(4, 2)
$rng = Get-CellDocRange $cell
$rngBorro = Split-RunAt $rng "borro"
$rngBorro.Text = "borró"

# Fila 5, Col 2: "Se borro la solicitud de la lista" -> "borró" y el
# bookmark _GoBack se mueve a continuacion de la palabra corregida
$cell = $t.Cell(5, 2)
$rng = Get-CellDocRange $cell
$rngBorro = Split-RunAt $rng "borro"
$rngBorro.Text = "borró"

# El bookmark _GoBack original esta al final de la ultima celda "OK"
# (fila 5, col 3); Word lo recoloca automaticamente en la ultima posicion
# editada, que ahora es justo despues de "borró" en la fila 5.
$okCell = $t.Cell(5, 3)
$okStart = $okCell.Range.Start
$okEnd = $okCell.Range.End
$okRng = $d.Range($okStart, $okEnd)
$okRng.Find.Execute("_GoBack_MARKER_NEVER_MATCHES", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$rngBorroEnd = $rngBorro.Duplicate
$rngBorroEnd.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngBorroEnd) | Out-Null
